$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.88
$ws.Range("G3").Value = 10.69
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "0.1284"
$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "1.429511"
$ws.Range("G4").Value = 13.37
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "0.1392"
$ws.Range("G5").Value = 12.5
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "0.1318"
$ws.Range("G6").Value = 9.970000000000001
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "0.1453"
$ws.Range("G7").Value = 6.83
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1.773467"
$ws.Range("G8").Value = 9.569800000000001
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = "0.3339"
$ws.Range("Q8").Value = ""
$ws.Range("G9").Value = 142.22
$ws.Range("G10").Value = 9.25
$ws.Range("G11").Value = 10.82
$ws.Range("M11").NumberFormat = "@"
$ws.Range("M11").Value = "0.1306"
$ws.Range("G12").Value = 19.51
$ws.Range("G13").Value = 5.09
$ws.Range("G14").Value = 9.1
$ws.Range("G15").Value = 8.48
$ws.Range("G16").Value = 5.06
$ws.Range("G17").Value = 9.75
$ws.Range("G18").Value = 8.18
$ws.Range("H18").Value = 114.52
$ws.Range("M18").NumberFormat = "@"
$ws.Range("M18").Value = "0.1515"
$ws.Range("G19").Value = 4.15
$ws.Range("M19").NumberFormat = "@"
$ws.Range("M19").Value = "0.1435"
$ws.Range("O19").Value = 4.09
$ws.Range("Q19").NumberFormat = "@"
$ws.Range("Q19").Value = "1.295808"
$ws.Range("G20").Value = 5.93
$ws.Range("Q20").NumberFormat = "@"
$ws.Range("Q20").Value = "1.623138"
$ws.Range("G21").Value = 10.99
$ws.Range("G22").Value = 11.25
$ws.Range("P22").Value = 18.88
$ws.Range("R22").NumberFormat = "@"
$ws.Range("R22").Value = "{'earningsDate': ['2023-11-07', '2023-11-11'], 'earningsAverage': 0.33, 'earningsLow': 0.14, 'earningsHigh': 0.43, 'revenueAverage': 77420000, 'revenueLow': 74000000, 'revenueHigh': 81770000}"
$ws.Range("G23").Value = 8.44
$ws.Range("M23").NumberFormat = "@"
$ws.Range("M23").Value = "0.168"
$ws.Range("G24").Value = 4.011
$ws.Range("G25").Value = 3.95
$ws.Range("G26").Value = 7.03
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "MONITOR - PRE"
$ws.Range("G27").Value = 10.44
$ws.Range("H27").Value = 156.6
$ws.Range("M27").NumberFormat = "@"
$ws.Range("M27").Value = "0.1633"
$ws.Range("G28").Value = 13.37
$ws.Range("M28").NumberFormat = "@"
$ws.Range("M28").Value = "0.13520001"
$ws.Range("R28").NumberFormat = "@"
$ws.Range("R28").Value = "{'earningsDate': ['2023-08-08'], 'earningsAverage': 0.44, 'earningsLow': 0.37, 'earningsHigh': 0.48, 'revenueAverage': 45620000, 'revenueLow': 26550000, 'revenueHigh': 75200000}"
$ws.Range("G29").Value = 45.48
$ws.Range("H29").Value = 45.48
$ws.Range("G30").Value = 5.86
$ws.Range("M30").NumberFormat = "@"
$ws.Range("M30").Value = "0.1402"
$ws.Range("G31").Value = 17.26
$ws.Range("G32").Value = 20.41
$ws.Range("M32").NumberFormat = "@"
$ws.Range("M32").Value = "0.1372"
$ws.Range("R32").NumberFormat = "@"
$ws.Range("R32").Value = "{'earningsDate': ['2023-08-08'], 'earningsAverage': 0.75, 'earningsLow': 0.75, 'earningsHigh': 0.78, 'revenueAverage': 447050000, 'revenueLow': 438000000, 'revenueHigh': 453250000}"
$ws.Range("G33").Value = 9.56
$ws.Range("M33").NumberFormat = "@"
$ws.Range("M33").Value = "0.15689999"
$ws.Range("P33").Value = 12.25
$ws.Range("G34").Value = 5.6
$ws.Range("G35").Value = 6.28
$ws.Range("M35").NumberFormat = "@"
$ws.Range("M35").Value = "0.1401"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "MONITOR - PRE"
$ws.Range("G36").Value = 14.01
$ws.Range("H36").Value = 70.05
$ws.Range("L36").NumberFormat = "@"
$ws.Range("L36").Value = "2023-08-15"
$ws.Range("M36").NumberFormat = "@"
$ws.Range("M36").Value = "0.112799995"
$ws.Range("R36").NumberFormat = "@"
$ws.Range("R36").Value = "{'earningsDate': ['2023-11-07', '2023-11-13'], 'earningsAverage': 0.24, 'earningsLow': 0.13, 'earningsHigh': 0.4, 'revenueAverage': 68290000, 'revenueLow': 58390000, 'revenueHigh': 98920000}"
$ws.Range("G37").Value = 11.12
$ws.Range("G38").Value = 7.71
$ws.Range("H38").Value = 1079.4
$ws.Range("M38").NumberFormat = "@"
$ws.Range("M38").Value = "0.1634"
$ws.Range("G39").Value = 131.53
$ws.Range("H39").Value = 526.12
$ws.Range("G40").Value = 5.63
$ws.Range("G41").Value = 8.33
$ws.Range("G42").Value = 6.55
$ws.Range("G43").Value = 3.9352
$ws.Range("G44").Value = 23.66
$ws.Range("M44").NumberFormat = "@"
$ws.Range("M44").Value = "0.1691"
$ws.Range("G45").Value = 11.65
$ws.Range("M45").NumberFormat = "@"
$ws.Range("M45").Value = "0.14"
$ws.Range("P45").Value = 17.73
$ws.Range("G46").Value = 1.33
$ws.Range("M46").NumberFormat = "@"
$ws.Range("M46").Value = "0.33080003"
$ws.Range("G47").Value = 25.85
$ws.Range("H47").Value = 1783.65
$ws.Range("G48").Value = 12.51
$ws.Range("G49").Value = 5.3
$ws.Range("M49").NumberFormat = "@"
$ws.Range("M49").Value = "0.1895"
$ws.Range("G50").Value = 6.47
$ws.Range("G51").Value = 1.61
$ws.Range("G52").Value = 10.83
$ws.Range("M52").NumberFormat = "@"
$ws.Range("M52").Value = "0.1246"
$ws.Range("P52").Value = 12.38
$ws.Range("Q52").NumberFormat = "@"
$ws.Range("Q52").Value = "1.980256"
$ws.Range("G53").Value = 18.5
$ws.Range("G54").Value = 4.34
$ws.Range("M54").NumberFormat = "@"
$ws.Range("M54").Value = "0.1382"
$ws.Range("G55").Value = 5.87
$ws.Range("G56").Value = 19.96
$ws.Range("M56").NumberFormat = "@"
$ws.Range("M56").Value = "0.1311"
$ws.Range("G57").Value = 9.08
$ws.Range("G58").Value = 7.92
$ws.Range("G59").Value = 8.1
$ws.Range("M59").NumberFormat = "@"
$ws.Range("M59").Value = "0.21530001"
$ws.Range("G60").Value = 10.2
$ws.Range("G61").Value = 5.54
$ws.Range("H61").Value = 554
$ws.Range("G62").Value = 3.02
$ws.Range("M62").NumberFormat = "@"
$ws.Range("M62").Value = "0.1368"
$ws.Range("G63").Value = 13.47
$ws.Range("M63").NumberFormat = "@"
$ws.Range("M63").Value = "0.42900002"
$ws.Range("G64").Value = 7.21
$ws.Range("M64").NumberFormat = "@"
$ws.Range("M64").Value = "0.0696"
$ws.Range("G65").Value = 20.7
$ws.Range("M65").NumberFormat = "@"
$ws.Range("M65").Value = "0.1333"
$ws.Range("G66").Value = 7.33
$ws.Range("G67").Value = 7.78
$ws.Range("M67").NumberFormat = "@"
$ws.Range("M67").Value = "0.0825"
$ws.Range("G68").Value = 3.4
$ws.Range("M68").NumberFormat = "@"
$ws.Range("M68").Value = "0.1494"
$ws.Range("P68").Value = 5.105
$ws.Range("G69").Value = 18.54
$ws.Range("H69").Value = 1371.96
$ws.Range("G70").Value = 10.56
$ws.Range("M70").NumberFormat = "@"
$ws.Range("M70").Value = "0.1242"
$ws.Range("G71").Value = 38.04
$ws.Range("G72").Value = 9.390000000000001
$ws.Range("G73").Value = 11.32
$ws.Range("G74").Value = 14.17
$ws.Range("M75").NumberFormat = "@"
$ws.Range("M75").Value = "0.1341"
$ws.Range("G76").Value = 251.45
$ws.Range("H76").Value = 251.45
$ws.Range("G77").Value = 13.71
$ws.Range("G78").Value = 10.13
$ws.Range("M78").NumberFormat = "@"
$ws.Range("M78").Value = "0.1442"
$ws.Range("P78").Value = 15.8
$ws.Range("R78").NumberFormat = "@"
$ws.Range("R78").Value = "{'earningsDate': ['2023-11-01', '2023-11-06'], 'earningsAverage': {}, 'earningsLow': {}, 'earningsHigh': {}, 'revenueAverage': 5330000, 'revenueLow': 5330000, 'revenueHigh': 5330000}"
$ws.Range("G79").Value = 14.5
